$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.865.68'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.900.21'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.44%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.40'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.75%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.62'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.88%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.897.22'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.92'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.89%  '
$ws.Range('E11').Value = '  -3.26%  '
$ws.Range('E12').Value = '  -1.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000230'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.16'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.68%  '
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.381.04'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.55%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.815.46'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.25%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.51'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.15%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.879.04'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '430.77'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.62%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.91'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.46%  '
$ws.Range('E22').Value = '  -1.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.87'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '78.91'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.99'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.34%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.05'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -9.68%  '
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('E28').Value = '  -3.94%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000110'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +8.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.51'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.01%  '
$ws.Range('E32').Value = '  -6.52%  '
$ws.Range('E34').Value = '  -2.73%  '
$ws.Range('E35').Value = '  -2.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.954'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.57%  '
$ws.Range('E37').Value = '  -3.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '48.82'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.69%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.83'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.12%  '
$ws.Range('E40').Value = '  -4.84%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.114'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.22%  '
$ws.Range('E42').Value = '  -2.66%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.31'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.47%  '
$ws.Range('E44').Value = '  -2.95%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.702.02'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.16%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0336'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.04%  '
$ws.Range('E47').Value = '  -2.82%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '345.34'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.55'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.45%  '
